$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 12).Value = 'stimuli/img_r10cu.png'
$ws.Cells.Item(2, 13).Value = 78.52380952380952
$ws.Cells.Item(2, 14).Value = 56.14285714285715
$ws.Cells.Item(2, 15).Value = 67.33333333333333
$ws.Cells.Item(2, 16).Value = 42
$ws.Cells.Item(2, 17).Value = 7
$ws.Cells.Item(2, 18).Value = 7
$ws.Cells.Item(2, 19).Value = 7
$ws.Cells.Item(2, 20).Value = 7
$ws.Cells.Item(2, 21).Value = 7
$ws.Cells.Item(2, 22).Value = 6
$ws.Cells.Item(3, 9).Value = 'target'
$ws.Cells.Item(3, 10).Value = 'old'
$ws.Cells.Item(3, 11).Value = 'j'
$ws.Cells.Item(3, 12).Value = 'stimuli/img_vnxft.png'
$ws.Cells.Item(3, 13).Value = 53.22727272727273
$ws.Cells.Item(3, 14).Value = 34.84090909090909
$ws.Cells.Item(3, 15).Value = 44.03409090909091
$ws.Cells.Item(3, 16).Value = 44
$ws.Cells.Item(3, 17).Value = 3
$ws.Cells.Item(3, 18).Value = 3
$ws.Cells.Item(3, 19).Value = 3
$ws.Cells.Item(3, 20).Value = 4
$ws.Cells.Item(3, 21).Value = 3
$ws.Cells.Item(3, 22).Value = 4
$ws.Cells.Item(4, 9).Value = $null
$ws.Cells.Item(4, 10).Value = 'new'
$ws.Cells.Item(4, 11).Value = 'f'
$ws.Cells.Item(4, 12).Value = 'stimuli/img_s2zoe.png'
$ws.Cells.Item(4, 13).Value = 64.71428571428571
$ws.Cells.Item(4, 14).Value = 44.90476190476191
$ws.Cells.Item(4, 15).Value = 54.80952380952381
$ws.Cells.Item(4, 16).Value = 42
$ws.Cells.Item(4, 17).Value = 5
$ws.Cells.Item(4, 18).Value = 5
$ws.Cells.Item(4, 19).Value = 5
$ws.Cells.Item(4, 20).Value = 5
$ws.Cells.Item(4, 21).Value = 5
$ws.Cells.Item(4, 22).Value = 5
$ws.Cells.Item(5, 12).Value = 'stimuli/img_rru0v.png'
$ws.Cells.Item(5, 13).Value = 56.45238095238095
$ws.Cells.Item(5, 14).Value = 39.42857142857143
$ws.Cells.Item(5, 15).Value = 47.94047619047619
$ws.Cells.Item(5, 17).Value = 4
$ws.Cells.Item(5, 18).Value = 4
$ws.Cells.Item(5, 19).Value = 4
$ws.Cells.Item(5, 20).Value = 4
$ws.Cells.Item(5, 21).Value = 4
$ws.Cells.Item(5, 22).Value = 4
$ws.Cells.Item(6, 12).Value = 'stimuli/img_6wne4.png'
$ws.Cells.Item(6, 13).Value = 25.16279069767442
$ws.Cells.Item(6, 14).Value = 15
$ws.Cells.Item(6, 15).Value = 20.08139534883721
$ws.Cells.Item(6, 16).Value = 43
$ws.Cells.Item(6, 17).Value = 1
$ws.Cells.Item(6, 18).Value = 1
$ws.Cells.Item(6, 19).Value = 1
$ws.Cells.Item(6, 20).Value = 1
$ws.Cells.Item(6, 21).Value = 2
$ws.Cells.Item(6, 22).Value = 1
$ws.Cells.Item(7, 12).Value = 'stimuli/img_9684y.png'
$ws.Cells.Item(7, 13).Value = 77.95744680851064
$ws.Cells.Item(7, 14).Value = 56.70212765957447
$ws.Cells.Item(7, 15).Value = 67.32978723404256
$ws.Cells.Item(7, 16).Value = 47
$ws.Cells.Item(7, 17).Value = 7
$ws.Cells.Item(7, 18).Value = 7
$ws.Cells.Item(7, 19).Value = 7
$ws.Cells.Item(7, 20).Value = 6
$ws.Cells.Item(7, 21).Value = 7
$ws.Cells.Item(7, 22).Value = 7
$ws.Cells.Item(8, 12).Value = 'stimuli/img_vg73h.png'
$ws.Cells.Item(8, 13).Value = 87.7
$ws.Cells.Item(8, 14).Value = 72.4
$ws.Cells.Item(8, 15).Value = 80.05000000000001
$ws.Cells.Item(8, 16).Value = 50
$ws.Cells.Item(8, 17).Value = 10
$ws.Cells.Item(8, 18).Value = 10
$ws.Cells.Item(8, 19).Value = 10
$ws.Cells.Item(8, 20).Value = 10
$ws.Cells.Item(8, 21).Value = 10
$ws.Cells.Item(8, 22).Value = 10
$ws.Cells.Item(9, 12).Value = 'stimuli/img_37hgm.png'
$ws.Cells.Item(9, 13).Value = 70.95454545454545
$ws.Cells.Item(9, 14).Value = 54.77272727272727
$ws.Cells.Item(9, 15).Value = 62.86363636363636
$ws.Cells.Item(9, 16).Value = 44
$ws.Cells.Item(9, 17).Value = 6
$ws.Cells.Item(9, 18).Value = 6
$ws.Cells.Item(9, 19).Value = 6
$ws.Cells.Item(9, 20).Value = 6
$ws.Cells.Item(9, 21).Value = 6
$ws.Cells.Item(9, 22).Value = 6
$ws.Cells.Item(10, 9).Value = 'target'
$ws.Cells.Item(10, 10).Value = 'old'
$ws.Cells.Item(10, 11).Value = 'j'
$ws.Cells.Item(10, 12).Value = 'stimuli/img_bg264.png'
$ws.Cells.Item(10, 13).Value = 87.9047619047619
$ws.Cells.Item(10, 14).Value = 71.5
$ws.Cells.Item(10, 15).Value = 79.70238095238095
$ws.Cells.Item(10, 16).Value = 42
$ws.Cells.Item(10, 17).Value = 10
$ws.Cells.Item(10, 18).Value = 10
$ws.Cells.Item(10, 19).Value = 10
$ws.Cells.Item(10, 20).Value = 9
$ws.Cells.Item(10, 21).Value = 10
$ws.Cells.Item(10, 22).Value = 9
$ws.Cells.Item(12, 12).Value = 'stimuli/img_tbs4n.png'
$ws.Cells.Item(12, 13).Value = 78.95744680851064
$ws.Cells.Item(12, 14).Value = 58.97872340425532
$ws.Cells.Item(12, 15).Value = 68.96808510638297
$ws.Cells.Item(12, 16).Value = 47
$ws.Cells.Item(12, 17).Value = 7
$ws.Cells.Item(12, 18).Value = 7
$ws.Cells.Item(12, 19).Value = 7
$ws.Cells.Item(12, 20).Value = 7
$ws.Cells.Item(12, 21).Value = 7
$ws.Cells.Item(12, 22).Value = 7
$ws.Cells.Item(14, 9).Value = $null
$ws.Cells.Item(14, 10).Value = 'new'
$ws.Cells.Item(14, 11).Value = 'f'
$ws.Cells.Item(14, 12).Value = 'stimuli/img_5tr4v.png'
$ws.Cells.Item(14, 13).Value = 56.86046511627907
$ws.Cells.Item(14, 14).Value = 39.3953488372093
$ws.Cells.Item(14, 15).Value = 48.12790697674419
$ws.Cells.Item(14, 16).Value = 43
$ws.Cells.Item(14, 17).Value = 4
$ws.Cells.Item(14, 18).Value = 4
$ws.Cells.Item(14, 19).Value = 4
$ws.Cells.Item(14, 20).Value = 4
$ws.Cells.Item(14, 21).Value = 4
$ws.Cells.Item(14, 22).Value = 4
$ws.Cells.Item(15, 9).Value = 'target'
$ws.Cells.Item(15, 10).Value = 'old'
$ws.Cells.Item(15, 11).Value = 'j'
$ws.Cells.Item(15, 12).Value = 'stimuli/img_bf8nx.png'
$ws.Cells.Item(15, 13).Value = 86.63414634146342
$ws.Cells.Item(15, 14).Value = 66.63414634146342
$ws.Cells.Item(15, 15).Value = 76.63414634146342
$ws.Cells.Item(15, 16).Value = 41
$ws.Cells.Item(15, 17).Value = 9
$ws.Cells.Item(15, 18).Value = 9
$ws.Cells.Item(15, 19).Value = 9
$ws.Cells.Item(15, 20).Value = 8
$ws.Cells.Item(15, 21).Value = 9
$ws.Cells.Item(15, 22).Value = 8
$ws.Cells.Item(16, 8).Value = 'living_rooms'
$ws.Cells.Item(16, 9).Value = 'target'
$ws.Cells.Item(16, 10).Value = 'old'
$ws.Cells.Item(16, 11).Value = 'j'
$ws.Cells.Item(16, 12).Value = 'stimuli/img_qdln8.png'
$ws.Cells.Item(16, 13).Value = 85.51162790697674
$ws.Cells.Item(16, 14).Value = 67.86046511627907
$ws.Cells.Item(16, 15).Value = 76.68604651162791
$ws.Cells.Item(16, 16).Value = 43
$ws.Cells.Item(16, 17).Value = 9
$ws.Cells.Item(16, 18).Value = 9
$ws.Cells.Item(16, 19).Value = 9
$ws.Cells.Item(16, 20).Value = 9
$ws.Cells.Item(16, 21).Value = 9
$ws.Cells.Item(16, 22).Value = 9
$ws.Cells.Item(17, 12).Value = 'stimuli/img_4o8l0.png'
$ws.Cells.Item(17, 13).Value = 46.02173913043478
$ws.Cells.Item(17, 14).Value = 31.45652173913043
$ws.Cells.Item(17, 15).Value = 38.73913043478261
$ws.Cells.Item(17, 16).Value = 46
$ws.Cells.Item(17, 17).Value = 3
$ws.Cells.Item(17, 18).Value = 3
$ws.Cells.Item(17, 19).Value = 3
$ws.Cells.Item(17, 20).Value = 3
$ws.Cells.Item(17, 21).Value = 3
$ws.Cells.Item(17, 22).Value = 3
$ws.Cells.Item(18, 12).Value = 'stimuli/img_bj99b.png'
$ws.Cells.Item(18, 13).Value = 82.79069767441861
$ws.Cells.Item(18, 14).Value = 65.46511627906976
$ws.Cells.Item(18, 15).Value = 74.12790697674419
$ws.Cells.Item(18, 16).Value = 43
$ws.Cells.Item(18, 17).Value = 8
$ws.Cells.Item(18, 18).Value = 8
$ws.Cells.Item(18, 19).Value = 8
$ws.Cells.Item(18, 20).Value = 8
$ws.Cells.Item(18, 21).Value = 8
$ws.Cells.Item(18, 22).Value = 8
$ws.Cells.Item(20, 9).Value = $null
$ws.Cells.Item(20, 10).Value = 'new'
$ws.Cells.Item(20, 11).Value = 'f'
$ws.Cells.Item(20, 12).Value = 'stimuli/img_bbs77.png'
$ws.Cells.Item(20, 13).Value = 31.64444444444445
$ws.Cells.Item(20, 14).Value = 21.26666666666667
$ws.Cells.Item(20, 15).Value = 26.45555555555556
$ws.Cells.Item(20, 16).Value = 45
$ws.Cells.Item(20, 17).Value = 2
$ws.Cells.Item(20, 18).Value = 2
$ws.Cells.Item(20, 19).Value = 2
$ws.Cells.Item(20, 20).Value = 2
$ws.Cells.Item(20, 21).Value = 2
$ws.Cells.Item(20, 22).Value = 2
$ws.Cells.Item(21, 9).Value = 'target'
$ws.Cells.Item(21, 10).Value = 'old'
$ws.Cells.Item(21, 11).Value = 'j'
$ws.Cells.Item(21, 12).Value = 'stimuli/img_di6f0.png'
$ws.Cells.Item(21, 13).Value = 94.04347826086956
$ws.Cells.Item(21, 14).Value = 83.34782608695652
$ws.Cells.Item(21, 15).Value = 88.69565217391303
$ws.Cells.Item(21, 16).Value = 46
$ws.Cells.Item(21, 17).Value = 10
$ws.Cells.Item(21, 18).Value = 10
$ws.Cells.Item(21, 19).Value = 10
$ws.Cells.Item(21, 20).Value = 10
$ws.Cells.Item(21, 21).Value = 10
$ws.Cells.Item(21, 22).Value = 10
$ws.Cells.Item(22, 12).Value = 'stimuli/img_2qhro.png'
$ws.Cells.Item(22, 13).Value = 81.73809523809524
$ws.Cells.Item(22, 14).Value = 62.73809523809524
$ws.Cells.Item(22, 15).Value = 72.23809523809524
$ws.Cells.Item(22, 16).Value = 42
$ws.Cells.Item(22, 17).Value = 8
$ws.Cells.Item(22, 18).Value = 8
$ws.Cells.Item(22, 19).Value = 8
$ws.Cells.Item(22, 20).Value = 8
$ws.Cells.Item(22, 21).Value = 8
$ws.Cells.Item(22, 22).Value = 8
$ws.Cells.Item(23, 8).Value = $null
$ws.Cells.Item(23, 9).Value = $null
$ws.Cells.Item(23, 10).Value = 'catch'
$ws.Cells.Item(23, 11).Value = 'f'
$ws.Cells.Item(23, 12).Value = 'stimuli/catch_08.jpg'
$ws.Cells.Item(23, 13).Value = $null
$ws.Cells.Item(23, 14).Value = $null
$ws.Cells.Item(23, 15).Value = $null
$ws.Cells.Item(23, 16).Value = $null
$ws.Cells.Item(23, 17).Value = $null
$ws.Cells.Item(23, 18).Value = $null
$ws.Cells.Item(23, 19).Value = $null
$ws.Cells.Item(23, 20).Value = $null
$ws.Cells.Item(23, 21).Value = $null
$ws.Cells.Item(23, 22).Value = $null
$ws.Cells.Item(24, 9).Value = $null
$ws.Cells.Item(24, 10).Value = 'new'
$ws.Cells.Item(24, 11).Value = 'f'
$ws.Cells.Item(24, 12).Value = 'stimuli/img_lzz3x.png'
$ws.Cells.Item(24, 13).Value = 18.46341463414634
$ws.Cells.Item(24, 14).Value = 11.92682926829268
$ws.Cells.Item(24, 15).Value = 15.19512195121951
$ws.Cells.Item(24, 16).Value = 41
$ws.Cells.Item(24, 17).Value = 1
$ws.Cells.Item(24, 18).Value = 1
$ws.Cells.Item(24, 19).Value = 1
$ws.Cells.Item(24, 20).Value = 1
$ws.Cells.Item(24, 21).Value = 1
$ws.Cells.Item(24, 22).Value = 1
$ws.Cells.Item(25, 9).Value = $null
$ws.Cells.Item(25, 10).Value = 'new'
$ws.Cells.Item(25, 11).Value = 'f'
$ws.Cells.Item(25, 12).Value = 'stimuli/img_iudc4.png'
$ws.Cells.Item(25, 13).Value = 73.625
$ws.Cells.Item(25, 14).Value = 52.275
$ws.Cells.Item(25, 15).Value = 62.95
$ws.Cells.Item(25, 16).Value = 40
$ws.Cells.Item(25, 17).Value = 6
$ws.Cells.Item(25, 18).Value = 6
$ws.Cells.Item(25, 19).Value = 6
$ws.Cells.Item(25, 20).Value = 6
$ws.Cells.Item(25, 21).Value = 6
$ws.Cells.Item(25, 22).Value = 6
$ws.Cells.Item(26, 12).Value = 'stimuli/img_5nlnv.png'
$ws.Cells.Item(26, 13).Value = 86.1219512195122
$ws.Cells.Item(26, 14).Value = 69.1951219512195
$ws.Cells.Item(26, 15).Value = 77.65853658536585
$ws.Cells.Item(26, 16).Value = 41
$ws.Cells.Item(26, 17).Value = 9
$ws.Cells.Item(26, 18).Value = 9
$ws.Cells.Item(26, 19).Value = 9
$ws.Cells.Item(26, 20).Value = 9
$ws.Cells.Item(26, 21).Value = 9
$ws.Cells.Item(26, 22).Value = 9
$ws.Cells.Item(27, 12).Value = 'stimuli/img_il020.png'
$ws.Cells.Item(27, 13).Value = 18.85416666666667
$ws.Cells.Item(27, 14).Value = 16.16666666666667
$ws.Cells.Item(27, 15).Value = 17.51041666666667
$ws.Cells.Item(27, 16).Value = 48
$ws.Cells.Item(27, 17).Value = 1
$ws.Cells.Item(27, 18).Value = 1
$ws.Cells.Item(27, 19).Value = 1
$ws.Cells.Item(27, 20).Value = 1
$ws.Cells.Item(27, 21).Value = 1
$ws.Cells.Item(27, 22).Value = 1
$ws.Cells.Item(28, 12).Value = 'stimuli/img_eiu3c.png'
$ws.Cells.Item(28, 13).Value = 65.1590909090909
$ws.Cells.Item(28, 14).Value = 46.22727272727273
$ws.Cells.Item(28, 15).Value = 55.69318181818181
$ws.Cells.Item(28, 16).Value = 44
$ws.Cells.Item(29, 12).Value = 'stimuli/img_kq9s9.png'
$ws.Cells.Item(29, 13).Value = 62.30232558139535
$ws.Cells.Item(29, 14).Value = 39.97674418604651
$ws.Cells.Item(29, 15).Value = 51.13953488372093
$ws.Cells.Item(29, 20).Value = 5
$ws.Cells.Item(29, 21).Value = 5
$ws.Cells.Item(30, 9).Value = $null
$ws.Cells.Item(30, 10).Value = 'new'
$ws.Cells.Item(30, 11).Value = 'f'
$ws.Cells.Item(30, 12).Value = 'stimuli/img_jpjeg.png'
$ws.Cells.Item(30, 13).Value = 90.90697674418605
$ws.Cells.Item(30, 14).Value = 74.3953488372093
$ws.Cells.Item(30, 15).Value = 82.65116279069767
$ws.Cells.Item(30, 16).Value = 43
$ws.Cells.Item(31, 12).Value = 'stimuli/img_9oofc.png'
$ws.Cells.Item(31, 13).Value = 82.47619047619048
$ws.Cells.Item(31, 14).Value = 65.5
$ws.Cells.Item(31, 15).Value = 73.98809523809524
$ws.Cells.Item(31, 16).Value = 42
$ws.Cells.Item(32, 12).Value = 'stimuli/img_njhlh.png'
$ws.Cells.Item(32, 13).Value = 59.74418604651163
$ws.Cells.Item(32, 14).Value = 41.51162790697674
$ws.Cells.Item(32, 15).Value = 50.62790697674419
$ws.Cells.Item(32, 16).Value = 43
$ws.Cells.Item(32, 17).Value = 4
$ws.Cells.Item(32, 18).Value = 4
$ws.Cells.Item(32, 19).Value = 4
$ws.Cells.Item(32, 20).Value = 4
$ws.Cells.Item(32, 21).Value = 4
$ws.Cells.Item(32, 22).Value = 4
$ws.Cells.Item(33, 9).Value = 'target'
$ws.Cells.Item(33, 10).Value = 'old'
$ws.Cells.Item(33, 11).Value = 'j'
$ws.Cells.Item(33, 12).Value = 'stimuli/img_sx68r.png'
$ws.Cells.Item(33, 13).Value = 54
$ws.Cells.Item(33, 14).Value = 33.2051282051282
$ws.Cells.Item(33, 15).Value = 43.6025641025641
$ws.Cells.Item(33, 16).Value = 39
$ws.Cells.Item(33, 17).Value = 3
$ws.Cells.Item(33, 18).Value = 3
$ws.Cells.Item(33, 19).Value = 3
$ws.Cells.Item(33, 20).Value = 3
$ws.Cells.Item(33, 21).Value = 4
$ws.Cells.Item(33, 22).Value = 3
$ws.Cells.Item(34, 12).Value = 'stimuli/img_x9w7o.png'
$ws.Cells.Item(34, 13).Value = 92.38888888888889
$ws.Cells.Item(34, 14).Value = 72.94444444444444
$ws.Cells.Item(34, 15).Value = 82.66666666666666
$ws.Cells.Item(34, 16).Value = 36
$ws.Cells.Item(34, 17).Value = 10
$ws.Cells.Item(34, 18).Value = 10
$ws.Cells.Item(34, 19).Value = 10
$ws.Cells.Item(34, 20).Value = 10
$ws.Cells.Item(34, 21).Value = 10
$ws.Cells.Item(34, 22).Value = 10
$ws.Cells.Item(35, 9).Value = $null
$ws.Cells.Item(35, 10).Value = 'new'
$ws.Cells.Item(35, 11).Value = 'f'
$ws.Cells.Item(35, 12).Value = 'stimuli/img_qz292.png'
$ws.Cells.Item(35, 13).Value = 78.26666666666667
$ws.Cells.Item(35, 14).Value = 59.13333333333333
$ws.Cells.Item(35, 15).Value = 68.7
$ws.Cells.Item(35, 16).Value = 45
$ws.Cells.Item(35, 17).Value = 7
$ws.Cells.Item(35, 18).Value = 7
$ws.Cells.Item(35, 19).Value = 7
$ws.Cells.Item(35, 20).Value = 7
$ws.Cells.Item(35, 21).Value = 7
$ws.Cells.Item(35, 22).Value = 7
$ws.Cells.Item(36, 12).Value = 'stimuli/img_196rk.png'
$ws.Cells.Item(36, 13).Value = 86.53488372093024
$ws.Cells.Item(36, 14).Value = 69.46511627906976
$ws.Cells.Item(36, 15).Value = 78
$ws.Cells.Item(36, 17).Value = 9
$ws.Cells.Item(36, 18).Value = 9
$ws.Cells.Item(36, 19).Value = 9
$ws.Cells.Item(36, 20).Value = 9
$ws.Cells.Item(36, 21).Value = 9
$ws.Cells.Item(36, 22).Value = 9
$ws.Cells.Item(37, 9).Value = 'target'
$ws.Cells.Item(37, 10).Value = 'old'
$ws.Cells.Item(37, 11).Value = 'j'
$ws.Cells.Item(37, 12).Value = 'stimuli/img_tujn3.png'
$ws.Cells.Item(37, 13).Value = 81.4090909090909
$ws.Cells.Item(37, 14).Value = 62.52272727272727
$ws.Cells.Item(37, 15).Value = 71.9659090909091
$ws.Cells.Item(37, 16).Value = 44
$ws.Cells.Item(37, 17).Value = 8
$ws.Cells.Item(37, 18).Value = 8
$ws.Cells.Item(37, 19).Value = 8
$ws.Cells.Item(37, 20).Value = 7
$ws.Cells.Item(37, 21).Value = 8
$ws.Cells.Item(37, 22).Value = 7
$ws.Cells.Item(38, 9).Value = $null
$ws.Cells.Item(38, 10).Value = 'new'
$ws.Cells.Item(38, 11).Value = 'f'
$ws.Cells.Item(38, 12).Value = 'stimuli/img_rg4in.png'
$ws.Cells.Item(38, 13).Value = 49.3695652173913
$ws.Cells.Item(38, 14).Value = 30.21739130434782
$ws.Cells.Item(38, 15).Value = 39.79347826086956
$ws.Cells.Item(38, 16).Value = 46
$ws.Cells.Item(38, 17).Value = 3
$ws.Cells.Item(38, 18).Value = 3
$ws.Cells.Item(38, 19).Value = 3
$ws.Cells.Item(38, 20).Value = 3
$ws.Cells.Item(38, 21).Value = 3
$ws.Cells.Item(38, 22).Value = 3
$ws.Cells.Item(39, 12).Value = 'stimuli/img_99exx.png'
$ws.Cells.Item(39, 13).Value = 70.02272727272727
$ws.Cells.Item(39, 14).Value = 51.88636363636363
$ws.Cells.Item(39, 15).Value = 60.95454545454545
$ws.Cells.Item(39, 17).Value = 6
$ws.Cells.Item(39, 18).Value = 6
$ws.Cells.Item(39, 19).Value = 6
$ws.Cells.Item(39, 20).Value = 5
$ws.Cells.Item(39, 21).Value = 5
$ws.Cells.Item(39, 22).Value = 6
$ws.Cells.Item(40, 9).Value = 'target'
$ws.Cells.Item(40, 10).Value = 'old'
$ws.Cells.Item(40, 11).Value = 'j'
$ws.Cells.Item(40, 12).Value = 'stimuli/img_j73b6.png'
$ws.Cells.Item(40, 13).Value = 21.5609756097561
$ws.Cells.Item(40, 14).Value = 19.90243902439024
$ws.Cells.Item(40, 15).Value = 20.73170731707317
$ws.Cells.Item(40, 17).Value = 1
$ws.Cells.Item(40, 18).Value = 1
$ws.Cells.Item(40, 19).Value = 1
$ws.Cells.Item(40, 20).Value = 2
$ws.Cells.Item(40, 21).Value = 1
$ws.Cells.Item(40, 22).Value = 2
$ws.Cells.Item(41, 9).Value = 'target'
$ws.Cells.Item(41, 10).Value = 'old'
$ws.Cells.Item(41, 11).Value = 'j'
$ws.Cells.Item(41, 12).Value = 'stimuli/img_xbtev.png'
$ws.Cells.Item(41, 13).Value = 13.68181818181818
$ws.Cells.Item(41, 14).Value = 8.568181818181818
$ws.Cells.Item(41, 15).Value = 11.125
$ws.Cells.Item(41, 16).Value = 44
$ws.Cells.Item(41, 17).Value = 1
$ws.Cells.Item(41, 18).Value = 1
$ws.Cells.Item(41, 19).Value = 1
$ws.Cells.Item(41, 20).Value = 1
$ws.Cells.Item(41, 21).Value = 1
$ws.Cells.Item(41, 22).Value = 1
$ws.Cells.Item(42, 12).Value = 'stimuli/img_pey7u.png'
$ws.Cells.Item(42, 13).Value = 30.34883720930232
$ws.Cells.Item(42, 14).Value = 20.34883720930232
$ws.Cells.Item(42, 15).Value = 25.34883720930232
$ws.Cells.Item(42, 16).Value = 43
$ws.Cells.Item(42, 17).Value = 1
$ws.Cells.Item(42, 18).Value = 2
$ws.Cells.Item(42, 19).Value = 2
$ws.Cells.Item(42, 20).Value = 2
$ws.Cells.Item(42, 21).Value = 2
$ws.Cells.Item(42, 22).Value = 2
